$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 33 (Sarango Hidalgo Pablo Fernando) ---
$ws.Range("L33").Value = 1
$ws.Range("Q33").Value = 6
$ws.Range("R33").Value = 8.75820382759259

# --- Append new student row 37 (Ana Salet Hidalgo Flores) ---
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "Ana Salet Hidalgo Flores"
$ws.Range("C37").Value = "anitasalet2203@gmail.com"
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 1
$ws.Range("I37").Value = 3
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 1
$ws.Range("N37").Value = 0.8807970779778823
$ws.Range("O37").Value = 1
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0
$ws.Range("R37").Value = 2.880797077977882
